$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrower status columns (width 17.216 -> 13.410 "characters") ----------
# ColumnWidth is snapped to whole pixels by Excel, 12.5 is the input that
# lands on the closest attainable snapped width to 13.4101845877511.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
